# Game/Assets/Localization/LanguageData.xlsx
# Fill in the localized text for the three new missions (fireball tutorial,
# boomerang, and the final boss) that were previously placeholder/empty
# cells in columns B (English) and C (Portuguese) for rows 35-43, and
# restore the view state (scroll position / active cell) that the author
# left the workbook in after testing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mission 03: fireball tutorial -----------------------------------
$ws.Range("B35").Value = "Learn how to use the fireball!"
$ws.Range("C35").Value = "Aprenda a usar a Bola de Fogo!"

$ws.Range("B36").Value = "Get the fireball spell, and use it to explore the abandoned ruins in the forest."
$ws.Range("C36").Value = "Obtenha a magia da bola de fogo e use-a para explorar as ruínas na floresta."

$ws.Range("B37").Value = "You learned how to use the fireball!"
$ws.Range("C37").Value = "Você aprendeu a usar a bola de fogo!"

# --- Mission 04: get the boomerang -------------------------------------
$ws.Range("B38").Value = "Get the Boomerang"
$ws.Range("C38").Value = "Obtenha o bumerangue"

$ws.Range("B39").Value = "Get the boomerang on the dark forest, and learn how to control it to defeat enemies."
$ws.Range("C39").Value = "Obtenha o bumerangue na floresta negra, e aprenda a controlá-lo para destruir inimigos."

$ws.Range("B40").Value = "You got the boomerang!"
$ws.Range("C40").Value = "Você obteve o bumerangue!"

# --- Mission 05: kill the boss ------------------------------------------
$ws.Range("B41").Value = "Kill the boss"
$ws.Range("C41").Value = "Mate o chefe"

$ws.Range("B42").Value = "Use all your weapons and knowledge to kill the boss."
$ws.Range("C42").Value = "Use suas armas e conhecimento para matar o chefe."

$ws.Range("B43").Value = "Boss killed!"
$ws.Range("C43").Value = "Chefe morto!"

# --- Restore the view/selection state left by the author ---------------
$ws.Activate()
$ws.Range("B41").Select()
$excel.ActiveWindow.ScrollRow = 25
$excel.ActiveWindow.ScrollColumn = 1
